$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add a new "2022-Q4" sheet, right after "总计" (i.e. before the
#    existing "2022-Q2" sheet). We clone the "2022-Q2" sheet so the
#    new sheet inherits the same layout/column styling, then we
#    overwrite its name + data with the 2022-Q4 figures.
# ------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item(2)
$sheetQ2.Copy($sheetQ2)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The Q2 template has 8 data rows (rows 2-9); the Q4 data only has 7
# data rows (rows 2-8), so drop the extra trailing row.
$newSheet.Rows.Item(9).Delete()

# Fund holdings for 2022-Q4.
$q4Data = @(
    @("005561", "创金合信中证红利低波动指数A", "3.22", "94.10", "2.43", "0.0782", 4),
    @("005562", "创金合信中证红利低波动指数C", "2.76", "94.10", "2.43", "0.0671", 4),
    @("512890", "华泰柏瑞中证红利低波动ETF", "2.58", "99.39", "2.58", "0.0666", 4),
    @("009658", "汇丰晋信中小盘低波动策略股票A", "0.85", "92.42", "1.44", "0.0122", 8),
    @("006863", "国联安智能制造混合", "0.20", "94.39", "3.71", "0.0074", 9),
    @("005770", "信澳中证沪港深高股息精选指数", "0.13", "23.47", "0.70", "0.0009", 2),
    @("009775", "汇丰晋信中小盘低波动策略股票C", "0.04", "92.42", "1.44", "0.0006", 8)
)

# A plain, never-styled cell used as a format source: PasteSpecial-ing
# its (blank) format onto a cell strips the "quotePrefix" style that
# Excel otherwise stamps on any cell entered with a leading apostrophe,
# so the text cells below end up with no style index at all - matching
# freshly authored inline-string cells.
$blankStyleSrc = $newSheet.Range("C2")

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]
    $newSheet.Range("A$r").Value = $i
    $newSheet.Range("B$r").Value = "'" + $row[0]
    $newSheet.Range("C$r").Value = $row[1]
    $newSheet.Range("D$r").Value = "'" + $row[2]
    $newSheet.Range("E$r").Value = "'" + $row[3]
    $newSheet.Range("F$r").Value = "'" + $row[4]
    $newSheet.Range("G$r").Value = "'" + $row[5]
    $newSheet.Range("H$r").Value = $row[6]

    $blankStyleSrc.Copy()
    $newSheet.Range("B$r").PasteSpecial(-4122)
    $blankStyleSrc.Copy()
    $newSheet.Range("D$r`:G$r").PasteSpecial(-4122)
}

# ------------------------------------------------------------------
# 2) Insert the 2022-Q4 summary row at the top of the "总计" sheet's
#    data, pushing the existing quarters down by one row, and
#    renumber the running index column (A).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$lastRow = $total.UsedRange.Rows.Count

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.23

for ($r = 3; $r -le ($lastRow + 1); $r++) {
    $total.Range("A$r").Value = $r - 2
}

# Restore the original active sheet/selection so view-state metadata
# matches a normal "just opened the workbook" state.
$total.Range("A1").Select()
$total.Activate()
